# Refactor the "intervention_type" synthetic color array:
#   black square (noir)  -> blue book  (bleu)
#   red square            -> red book
#   orange square         -> orange book
#   green square          -> green book
#
# The word label "noir" becomes "bleu" while "rouge" / "orange" / "vert"
# stay as-is (only the black entry's word changes, same as the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "⬛"   = "📘"
    "🟥"   = "📕"
    "🟧"   = "📙"
    "🟩"   = "📗"
    "noir" = "bleu"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Text
        if ($map.ContainsKey($text)) {
            $cell.Value = $map[$text]
        }
    }
}
